# Updated readme and summary
# - Point B5 at the new repo path (jferrari's old machine -> julia's gitHub checkout)
# - Highlight the two "No module named X" rows (21-24) with a purple fill
# - Leave selection on the newly-highlighted block, matching the saved view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "C:\Users\julia\gitHub\CompuCell3D\CompuCell3D\core\"

$ws.Range("A21:C24").Interior.Color = 10498160

$ws.Range("A21:C24").Select() | Out-Null
